# Replace the numeric OrderNumber values in column A (rows 2-21) with
# placeholder text labels "delete01".."delete20" (these get written out as
# shared strings), then move the active selection to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 21; $row++) {
    $n = $row - 1
    $ws.Range("A$row").Value = "delete{0:D2}" -f $n
}

$ws.Range("C2").Select()
